$d = $word.ActiveDocument

# Locate, by index, the first paragraph to delete ("23,45,34,1,2") and the
# paragraph containing "4 - 3" (the last one is the trailing " " paragraph
# right after it, which carries the lastRenderedPageBreak and precedes the
# section break). We walk $d.Paragraphs by index rather than using
# Range.Paragraphs, since ranges derived from Find results do not reliably
# report paragraph text in this runtime.
$count = $d.Paragraphs.Count
$startIdx = -1
$fourThreeIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startIdx -eq -1 -and $t.Contains("23,45,34,1,2")) {
        $startIdx = $i
    }
    if ($t.Contains("4") -and $t.Contains([char]0x2013) -and $t.Contains("3")) {
        $fourThreeIdx = $i
    }
}
$endIdx = $fourThreeIdx + 1
$beforeIdx = $startIdx - 1

# Insert the two replacement blank paragraphs right after the paragraph that
# precedes the block being removed *before* deleting anything. (Inserting
# after the delete would make that anchor paragraph the very last paragraph
# in the document body, and a zero-length Range at the final paragraph mark
# inserts new content *before* it rather than after.)
$beforePara = $d.Paragraphs.Item($beforeIdx)
$insertPoint = $d.Range($beforePara.Range.End, $beforePara.Range.End)
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
[void]$insertPoint.InsertXML($blankParaXml + $blankParaXml)

# The two new blank paragraphs shifted everything that followed down by two
# indices, so re-resolve the start/end paragraphs of the block to delete.
$startPara = $d.Paragraphs.Item($startIdx + 2)
$endPara = $d.Paragraphs.Item($endIdx + 2)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
